# Auto-generated edit script: add 2022-Q3 sheet + summary row
$wb = $excel.ActiveWorkbook

# --- Step 1: insert new sheet "2022-Q3" before "2022-Q2" ---
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$newQ3 = $wb.Worksheets.Add($wsQ2)
$newQ3.Name = "2022-Q3"

# Copy header row + column-A formatting from "2022-Q2" (same layout family)
$wsQ2.Rows("1:1").Copy()
$newQ3.Range("A1:H1").PasteSpecial(-4122)
$wsQ2.Range("A2").Copy()
$newQ3.Range("A2:A40").PasteSpecial(-4122)

# Header row (all Chinese text, never numeric-looking)
$newQ3.Range("B1").Value = '基金代码'
$newQ3.Range("C1").Value = '基金名称'
$newQ3.Range("D1").Value = '基金规模'
$newQ3.Range("E1").Value = '股票总仓位'
$newQ3.Range("F1").Value = '仓位占比'
$newQ3.Range("G1").Value = '持有市值(亿元)'
$newQ3.Range("H1").Value = '仓位排名'

# Data rows
$newQ3.Range("A2").Value = 0
$c = $newQ3.Range("B2")
$c.NumberFormat = "@"
$c.Value = '016935'
$c = $newQ3.Range("D2")
$c.NumberFormat = "@"
$c.Value = '15.57'
$c = $newQ3.Range("E2")
$c.NumberFormat = "@"
$c.Value = '93.89'
$c = $newQ3.Range("F2")
$c.NumberFormat = "@"
$c.Value = '1.64'
$newQ3.Range("C2").Value = '景顺长城中证500指数增强C'
$c = $newQ3.Range("G2")
$c.NumberFormat = "@"
$c.Value = '0.2553'
$newQ3.Range("H2").Value = 10
$newQ3.Range("A3").Value = 1
$c = $newQ3.Range("B3")
$c.NumberFormat = "@"
$c.Value = '000978'
$c = $newQ3.Range("D3")
$c.NumberFormat = "@"
$c.Value = '7.14'
$c = $newQ3.Range("E3")
$c.NumberFormat = "@"
$c.Value = '93.64'
$c = $newQ3.Range("F3")
$c.NumberFormat = "@"
$c.Value = '1.86'
$newQ3.Range("C3").Value = '景顺长城量化精选股票'
$c = $newQ3.Range("G3")
$c.NumberFormat = "@"
$c.Value = '0.1328'
$newQ3.Range("H3").Value = 8
$newQ3.Range("A4").Value = 2
$c = $newQ3.Range("B4")
$c.NumberFormat = "@"
$c.Value = '561550'
$c = $newQ3.Range("D4")
$c.NumberFormat = "@"
$c.Value = '7.81'
$c = $newQ3.Range("E4")
$c.NumberFormat = "@"
$c.Value = '99.23'
$c = $newQ3.Range("F4")
$c.NumberFormat = "@"
$c.Value = '1.18'
$newQ3.Range("C4").Value = '华泰柏瑞中证500增强策略ETF'
$c = $newQ3.Range("G4")
$c.NumberFormat = "@"
$c.Value = '0.0922'
$newQ3.Range("H4").Value = 5
$newQ3.Range("A5").Value = 3
$c = $newQ3.Range("B5")
$c.NumberFormat = "@"
$c.Value = '014155'
$c = $newQ3.Range("D5")
$c.NumberFormat = "@"
$c.Value = '6.64'
$c = $newQ3.Range("E5")
$c.NumberFormat = "@"
$c.Value = '92.15'
$c = $newQ3.Range("F5")
$c.NumberFormat = "@"
$c.Value = '1.18'
$newQ3.Range("C5").Value = '国泰君安中证500指数增强A'
$c = $newQ3.Range("G5")
$c.NumberFormat = "@"
$c.Value = '0.0784'
$newQ3.Range("H5").Value = 9
$newQ3.Range("A6").Value = 4
$c = $newQ3.Range("B6")
$c.NumberFormat = "@"
$c.Value = '700001'
$c = $newQ3.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.82'
$c = $newQ3.Range("E6")
$c.NumberFormat = "@"
$c.Value = '91.01'
$c = $newQ3.Range("F6")
$c.NumberFormat = "@"
$c.Value = '3.73'
$newQ3.Range("C6").Value = '平安行业先锋混合'
$c = $newQ3.Range("G6")
$c.NumberFormat = "@"
$c.Value = '0.0679'
$newQ3.Range("H6").Value = 4
$newQ3.Range("A7").Value = 5
$c = $newQ3.Range("B7")
$c.NumberFormat = "@"
$c.Value = '012010'
$c = $newQ3.Range("D7")
$c.NumberFormat = "@"
$c.Value = '6.39'
$c = $newQ3.Range("E7")
$c.NumberFormat = "@"
$c.Value = '29.75'
$c = $newQ3.Range("F7")
$c.NumberFormat = "@"
$c.Value = '0.96'
$newQ3.Range("C7").Value = '富国泰享回报6个月持有期混合A'
$c = $newQ3.Range("G7")
$c.NumberFormat = "@"
$c.Value = '0.0613'
$newQ3.Range("H7").Value = 8
$newQ3.Range("A8").Value = 6
$c = $newQ3.Range("B8")
$c.NumberFormat = "@"
$c.Value = '014156'
$c = $newQ3.Range("D8")
$c.NumberFormat = "@"
$c.Value = '4.02'
$c = $newQ3.Range("E8")
$c.NumberFormat = "@"
$c.Value = '92.15'
$c = $newQ3.Range("F8")
$c.NumberFormat = "@"
$c.Value = '1.18'
$newQ3.Range("C8").Value = '国泰君安中证500指数增强C'
$c = $newQ3.Range("G8")
$c.NumberFormat = "@"
$c.Value = '0.0474'
$newQ3.Range("H8").Value = 9
$newQ3.Range("A9").Value = 7
$c = $newQ3.Range("B9")
$c.NumberFormat = "@"
$c.Value = '460009'
$c = $newQ3.Range("D9")
$c.NumberFormat = "@"
$c.Value = '4.22'
$c = $newQ3.Range("E9")
$c.NumberFormat = "@"
$c.Value = '93.06'
$c = $newQ3.Range("F9")
$c.NumberFormat = "@"
$c.Value = '1.07'
$newQ3.Range("C9").Value = '华泰柏瑞量化先行混合A'
$c = $newQ3.Range("G9")
$c.NumberFormat = "@"
$c.Value = '0.0452'
$newQ3.Range("H9").Value = 7
$newQ3.Range("A10").Value = 8
$c = $newQ3.Range("B10")
$c.NumberFormat = "@"
$c.Value = '008851'
$c = $newQ3.Range("D10")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c = $newQ3.Range("E10")
$c.NumberFormat = "@"
$c.Value = '64.77'
$c = $newQ3.Range("F10")
$c.NumberFormat = "@"
$c.Value = '1.37'
$newQ3.Range("C10").Value = '景顺长城量化对冲策略三个月定期开放灵活配置混合'
$c = $newQ3.Range("G10")
$c.NumberFormat = "@"
$c.Value = '0.0406'
$newQ3.Range("H10").Value = 4
$newQ3.Range("A11").Value = 9
$c = $newQ3.Range("B11")
$c.NumberFormat = "@"
$c.Value = '012879'
$c = $newQ3.Range("D11")
$c.NumberFormat = "@"
$c.Value = '3.33'
$c = $newQ3.Range("E11")
$c.NumberFormat = "@"
$c.Value = '90.73'
$c = $newQ3.Range("F11")
$c.NumberFormat = "@"
$c.Value = '1.11'
$newQ3.Range("C11").Value = '中信建投量化精选6个月持有期混合C'
$c = $newQ3.Range("G11")
$c.NumberFormat = "@"
$c.Value = '0.0370'
$newQ3.Range("H11").Value = 7
$newQ3.Range("A12").Value = 10
$c = $newQ3.Range("B12")
$c.NumberFormat = "@"
$c.Value = '009726'
$c = $newQ3.Range("D12")
$c.NumberFormat = "@"
$c.Value = '2.67'
$c = $newQ3.Range("E12")
$c.NumberFormat = "@"
$c.Value = '90.23'
$c = $newQ3.Range("F12")
$c.NumberFormat = "@"
$c.Value = '1.36'
$newQ3.Range("C12").Value = '招商中证500等权重指数增强A'
$c = $newQ3.Range("G12")
$c.NumberFormat = "@"
$c.Value = '0.0363'
$newQ3.Range("H12").Value = 10
$newQ3.Range("A13").Value = 11
$c = $newQ3.Range("B13")
$c.NumberFormat = "@"
$c.Value = '006441'
$c = $newQ3.Range("D13")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c = $newQ3.Range("E13")
$c.NumberFormat = "@"
$c.Value = '93.60'
$c = $newQ3.Range("F13")
$c.NumberFormat = "@"
$c.Value = '1.50'
$newQ3.Range("C13").Value = '中信建投中证500指数增强C'
$c = $newQ3.Range("G13")
$c.NumberFormat = "@"
$c.Value = '0.0338'
$newQ3.Range("H13").Value = 1
$newQ3.Range("A14").Value = 12
$c = $newQ3.Range("B14")
$c.NumberFormat = "@"
$c.Value = '006440'
$c = $newQ3.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.14'
$c = $newQ3.Range("E14")
$c.NumberFormat = "@"
$c.Value = '93.60'
$c = $newQ3.Range("F14")
$c.NumberFormat = "@"
$c.Value = '1.50'
$newQ3.Range("C14").Value = '中信建投中证500指数增强A'
$c = $newQ3.Range("G14")
$c.NumberFormat = "@"
$c.Value = '0.0321'
$newQ3.Range("H14").Value = 1
$newQ3.Range("A15").Value = 13
$c = $newQ3.Range("B15")
$c.NumberFormat = "@"
$c.Value = '013767'
$c = $newQ3.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.88'
$c = $newQ3.Range("E15")
$c.NumberFormat = "@"
$c.Value = '92.26'
$c = $newQ3.Range("F15")
$c.NumberFormat = "@"
$c.Value = '3.45'
$newQ3.Range("C15").Value = '平安价值回报混合A'
$c = $newQ3.Range("G15")
$c.NumberFormat = "@"
$c.Value = '0.0304'
$newQ3.Range("H15").Value = 4
$newQ3.Range("A16").Value = 14
$c = $newQ3.Range("B16")
$c.NumberFormat = "@"
$c.Value = '519618'
$c = $newQ3.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.02'
$c = $newQ3.Range("E16")
$c.NumberFormat = "@"
$c.Value = '21.05'
$c = $newQ3.Range("F16")
$c.NumberFormat = "@"
$c.Value = '0.95'
$newQ3.Range("C16").Value = '银河君信灵活配置混合I'
$c = $newQ3.Range("G16")
$c.NumberFormat = "@"
$c.Value = '0.0287'
$newQ3.Range("H16").Value = 6
$newQ3.Range("A17").Value = 15
$c = $newQ3.Range("B17")
$c.NumberFormat = "@"
$c.Value = '008422'
$c = $newQ3.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.69'
$c = $newQ3.Range("E17")
$c.NumberFormat = "@"
$c.Value = '61.89'
$c = $newQ3.Range("F17")
$c.NumberFormat = "@"
$c.Value = '3.29'
$newQ3.Range("C17").Value = '中融研发创新混合A'
$c = $newQ3.Range("G17")
$c.NumberFormat = "@"
$c.Value = '0.0227'
$newQ3.Range("H17").Value = 7
$newQ3.Range("A18").Value = 16
$c = $newQ3.Range("B18")
$c.NumberFormat = "@"
$c.Value = '012878'
$c = $newQ3.Range("D18")
$c.NumberFormat = "@"
$c.Value = '1.67'
$c = $newQ3.Range("E18")
$c.NumberFormat = "@"
$c.Value = '90.73'
$c = $newQ3.Range("F18")
$c.NumberFormat = "@"
$c.Value = '1.11'
$newQ3.Range("C18").Value = '中信建投量化精选6个月持有期混合A'
$c = $newQ3.Range("G18")
$c.NumberFormat = "@"
$c.Value = '0.0185'
$newQ3.Range("H18").Value = 7
$newQ3.Range("A19").Value = 17
$c = $newQ3.Range("B19")
$c.NumberFormat = "@"
$c.Value = '014305'
$c = $newQ3.Range("D19")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c = $newQ3.Range("E19")
$c.NumberFormat = "@"
$c.Value = '34.78'
$c = $newQ3.Range("F19")
$c.NumberFormat = "@"
$c.Value = '0.76'
$newQ3.Range("C19").Value = '华泰柏瑞中证500指数增强A'
$c = $newQ3.Range("G19")
$c.NumberFormat = "@"
$c.Value = '0.0167'
$newQ3.Range("H19").Value = 1
$newQ3.Range("A20").Value = 18
$c = $newQ3.Range("B20")
$c.NumberFormat = "@"
$c.Value = '010194'
$c = $newQ3.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.63'
$c = $newQ3.Range("E20")
$c.NumberFormat = "@"
$c.Value = '66.43'
$c = $newQ3.Range("F20")
$c.NumberFormat = "@"
$c.Value = '2.57'
$newQ3.Range("C20").Value = '博时睿祥15个月定期开放混合A'
$c = $newQ3.Range("G20")
$c.NumberFormat = "@"
$c.Value = '0.0162'
$newQ3.Range("H20").Value = 9
$newQ3.Range("A21").Value = 19
$c = $newQ3.Range("B21")
$c.NumberFormat = "@"
$c.Value = '009727'
$c = $newQ3.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.12'
$c = $newQ3.Range("E21")
$c.NumberFormat = "@"
$c.Value = '90.23'
$c = $newQ3.Range("F21")
$c.NumberFormat = "@"
$c.Value = '1.36'
$newQ3.Range("C21").Value = '招商中证500等权重指数增强C'
$c = $newQ3.Range("G21")
$c.NumberFormat = "@"
$c.Value = '0.0152'
$newQ3.Range("H21").Value = 10
$newQ3.Range("A22").Value = 20
$c = $newQ3.Range("B22")
$c.NumberFormat = "@"
$c.Value = '008423'
$c = $newQ3.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.40'
$c = $newQ3.Range("E22")
$c.NumberFormat = "@"
$c.Value = '61.89'
$c = $newQ3.Range("F22")
$c.NumberFormat = "@"
$c.Value = '3.29'
$newQ3.Range("C22").Value = '中融研发创新混合C'
$c = $newQ3.Range("G22")
$c.NumberFormat = "@"
$c.Value = '0.0132'
$newQ3.Range("H22").Value = 7
$newQ3.Range("A23").Value = 21
$c = $newQ3.Range("B23")
$c.NumberFormat = "@"
$c.Value = '010658'
$c = $newQ3.Range("D23")
$c.NumberFormat = "@"
$c.Value = '3.22'
$c = $newQ3.Range("E23")
$c.NumberFormat = "@"
$c.Value = '20.40'
$c = $newQ3.Range("F23")
$c.NumberFormat = "@"
$c.Value = '0.28'
$newQ3.Range("C23").Value = '海富通欣睿混合C'
$c = $newQ3.Range("G23")
$c.NumberFormat = "@"
$c.Value = '0.0090'
$newQ3.Range("H23").Value = 6
$newQ3.Range("A24").Value = 22
$c = $newQ3.Range("B24")
$c.NumberFormat = "@"
$c.Value = '010657'
$c = $newQ3.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.71'
$c = $newQ3.Range("E24")
$c.NumberFormat = "@"
$c.Value = '20.40'
$c = $newQ3.Range("F24")
$c.NumberFormat = "@"
$c.Value = '0.28'
$newQ3.Range("C24").Value = '海富通欣睿混合A'
$c = $newQ3.Range("G24")
$c.NumberFormat = "@"
$c.Value = '0.0076'
$newQ3.Range("H24").Value = 6
$newQ3.Range("A25").Value = 23
$c = $newQ3.Range("B25")
$c.NumberFormat = "@"
$c.Value = '011554'
$c = $newQ3.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.79'
$c = $newQ3.Range("E25")
$c.NumberFormat = "@"
$c.Value = '34.87'
$c = $newQ3.Range("F25")
$c.NumberFormat = "@"
$c.Value = '0.73'
$newQ3.Range("C25").Value = '海富通欣利混合A'
$c = $newQ3.Range("G25")
$c.NumberFormat = "@"
$c.Value = '0.0058'
$newQ3.Range("H25").Value = 5
$newQ3.Range("A26").Value = 24
$c = $newQ3.Range("B26")
$c.NumberFormat = "@"
$c.Value = '159620'
$c = $newQ3.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.36'
$c = $newQ3.Range("E26")
$c.NumberFormat = "@"
$c.Value = '91.92'
$c = $newQ3.Range("F26")
$c.NumberFormat = "@"
$c.Value = '1.42'
$newQ3.Range("C26").Value = '华夏中证智选500成长创新策略ETF'
$c = $newQ3.Range("G26")
$c.NumberFormat = "@"
$c.Value = '0.0051'
$newQ3.Range("H26").Value = 9
$newQ3.Range("A27").Value = 25
$c = $newQ3.Range("B27")
$c.NumberFormat = "@"
$c.Value = '519617'
$c = $newQ3.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.46'
$c = $newQ3.Range("E27")
$c.NumberFormat = "@"
$c.Value = '21.05'
$c = $newQ3.Range("F27")
$c.NumberFormat = "@"
$c.Value = '0.95'
$newQ3.Range("C27").Value = '银河君信灵活配置混合C'
$c = $newQ3.Range("G27")
$c.NumberFormat = "@"
$c.Value = '0.0044'
$newQ3.Range("H27").Value = 6
$newQ3.Range("A28").Value = 26
$c = $newQ3.Range("B28")
$c.NumberFormat = "@"
$c.Value = '001664'
$c = $newQ3.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.33'
$c = $newQ3.Range("E28")
$c.NumberFormat = "@"
$c.Value = '27.84'
$c = $newQ3.Range("F28")
$c.NumberFormat = "@"
$c.Value = '1.11'
$newQ3.Range("C28").Value = '平安鑫安混合A'
$c = $newQ3.Range("G28")
$c.NumberFormat = "@"
$c.Value = '0.0037'
$newQ3.Range("H28").Value = 5
$newQ3.Range("A29").Value = 27
$c = $newQ3.Range("B29")
$c.NumberFormat = "@"
$c.Value = '519616'
$c = $newQ3.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.36'
$c = $newQ3.Range("E29")
$c.NumberFormat = "@"
$c.Value = '21.05'
$c = $newQ3.Range("F29")
$c.NumberFormat = "@"
$c.Value = '0.95'
$newQ3.Range("C29").Value = '银河君信灵活配置混合A'
$c = $newQ3.Range("G29")
$c.NumberFormat = "@"
$c.Value = '0.0034'
$newQ3.Range("H29").Value = 6
$newQ3.Range("A30").Value = 28
$c = $newQ3.Range("B30")
$c.NumberFormat = "@"
$c.Value = '010246'
$c = $newQ3.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.25'
$c = $newQ3.Range("E30")
$c.NumberFormat = "@"
$c.Value = '93.06'
$c = $newQ3.Range("F30")
$c.NumberFormat = "@"
$c.Value = '1.07'
$newQ3.Range("C30").Value = '华泰柏瑞量化先行混合C'
$c = $newQ3.Range("G30")
$c.NumberFormat = "@"
$c.Value = '0.0027'
$newQ3.Range("H30").Value = 7
$newQ3.Range("A31").Value = 29
$c = $newQ3.Range("B31")
$c.NumberFormat = "@"
$c.Value = '006433'
$c = $newQ3.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.23'
$c = $newQ3.Range("E31")
$c.NumberFormat = "@"
$c.Value = '27.39'
$c = $newQ3.Range("F31")
$c.NumberFormat = "@"
$c.Value = '1.11'
$newQ3.Range("C31").Value = '平安鑫利灵活配置混合C'
$c = $newQ3.Range("G31")
$c.NumberFormat = "@"
$c.Value = '0.0026'
$newQ3.Range("H31").Value = 5
$newQ3.Range("A32").Value = 30
$c = $newQ3.Range("B32")
$c.NumberFormat = "@"
$c.Value = '011555'
$c = $newQ3.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.33'
$c = $newQ3.Range("E32")
$c.NumberFormat = "@"
$c.Value = '34.87'
$c = $newQ3.Range("F32")
$c.NumberFormat = "@"
$c.Value = '0.73'
$newQ3.Range("C32").Value = '海富通欣利混合C'
$c = $newQ3.Range("G32")
$c.NumberFormat = "@"
$c.Value = '0.0024'
$newQ3.Range("H32").Value = 5
$newQ3.Range("A33").Value = 31
$c = $newQ3.Range("B33")
$c.NumberFormat = "@"
$c.Value = '007049'
$c = $newQ3.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.16'
$c = $newQ3.Range("E33")
$c.NumberFormat = "@"
$c.Value = '27.84'
$c = $newQ3.Range("F33")
$c.NumberFormat = "@"
$c.Value = '1.11'
$newQ3.Range("C33").Value = '平安鑫安混合E'
$c = $newQ3.Range("G33")
$c.NumberFormat = "@"
$c.Value = '0.0018'
$newQ3.Range("H33").Value = 5
$newQ3.Range("A34").Value = 32
$c = $newQ3.Range("B34")
$c.NumberFormat = "@"
$c.Value = '003626'
$c = $newQ3.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.15'
$c = $newQ3.Range("E34")
$c.NumberFormat = "@"
$c.Value = '27.39'
$c = $newQ3.Range("F34")
$c.NumberFormat = "@"
$c.Value = '1.11'
$newQ3.Range("C34").Value = '平安鑫利灵活配置混合A'
$c = $newQ3.Range("G34")
$c.NumberFormat = "@"
$c.Value = '0.0017'
$newQ3.Range("H34").Value = 5
$newQ3.Range("A35").Value = 33
$c = $newQ3.Range("B35")
$c.NumberFormat = "@"
$c.Value = '013768'
$c = $newQ3.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.04'
$c = $newQ3.Range("E35")
$c.NumberFormat = "@"
$c.Value = '92.26'
$c = $newQ3.Range("F35")
$c.NumberFormat = "@"
$c.Value = '3.45'
$newQ3.Range("C35").Value = '平安价值回报混合C'
$c = $newQ3.Range("G35")
$c.NumberFormat = "@"
$c.Value = '0.0014'
$newQ3.Range("H35").Value = 4
$newQ3.Range("A36").Value = 34
$c = $newQ3.Range("B36")
$c.NumberFormat = "@"
$c.Value = '012011'
$c = $newQ3.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.05'
$c = $newQ3.Range("E36")
$c.NumberFormat = "@"
$c.Value = '29.75'
$c = $newQ3.Range("F36")
$c.NumberFormat = "@"
$c.Value = '0.96'
$newQ3.Range("C36").Value = '富国泰享回报6个月持有期混合C'
$c = $newQ3.Range("G36")
$c.NumberFormat = "@"
$c.Value = '0.0005'
$newQ3.Range("H36").Value = 8
$newQ3.Range("A37").Value = 35
$c = $newQ3.Range("B37")
$c.NumberFormat = "@"
$c.Value = '014306'
$c = $newQ3.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.07'
$c = $newQ3.Range("E37")
$c.NumberFormat = "@"
$c.Value = '34.78'
$c = $newQ3.Range("F37")
$c.NumberFormat = "@"
$c.Value = '0.76'
$newQ3.Range("C37").Value = '华泰柏瑞中证500指数增强C'
$c = $newQ3.Range("G37")
$c.NumberFormat = "@"
$c.Value = '0.0005'
$newQ3.Range("H37").Value = 1
$newQ3.Range("A38").Value = 36
$c = $newQ3.Range("B38")
$c.NumberFormat = "@"
$c.Value = '001665'
$c = $newQ3.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01'
$c = $newQ3.Range("E38")
$c.NumberFormat = "@"
$c.Value = '27.84'
$c = $newQ3.Range("F38")
$c.NumberFormat = "@"
$c.Value = '1.11'
$newQ3.Range("C38").Value = '平安鑫安混合C'
$c = $newQ3.Range("G38")
$c.NumberFormat = "@"
$c.Value = '0.0001'
$newQ3.Range("H38").Value = 5
$newQ3.Range("A39").Value = 37
$c = $newQ3.Range("B39")
$c.NumberFormat = "@"
$c.Value = '010195'
$c = $newQ3.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.00'
$c = $newQ3.Range("E39")
$c.NumberFormat = "@"
$c.Value = '66.43'
$c = $newQ3.Range("F39")
$c.NumberFormat = "@"
$c.Value = '2.57'
$newQ3.Range("C39").Value = '博时睿祥15个月定期开放混合C'
$newQ3.Range("G39").Value = 0
$newQ3.Range("H39").Value = 9
$newQ3.Range("A40").Value = 38
$c = $newQ3.Range("B40")
$c.NumberFormat = "@"
$c.Value = '006682'
$c = $newQ3.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.00'
$c = $newQ3.Range("E40")
$c.NumberFormat = "@"
$c.Value = '93.89'
$c = $newQ3.Range("F40")
$c.NumberFormat = "@"
$c.Value = '1.64'
$newQ3.Range("C40").Value = '景顺长城中证500指数增强A'
$newQ3.Range("G40").Value = 0
$newQ3.Range("H40").Value = 10

# --- Step 2: update "总计" summary sheet ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows("2:2").Insert()
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 39
$wsTotal.Range("D2").Value = 1.17

